# Plantilla_MantPriorVias.xlsx update
#  - "Priorización" becomes the active/selected tab (was "Presupuesto").
#  - Hoja1 (the hidden helper sheet backing the URCI list) gets a new
#    leading 0 value so the list runs 0.0 .. 1.0 instead of 0.1 .. 1.0.

$wb = $excel.ActiveWorkbook

$wsPresupuesto  = $wb.Worksheets.Item("Presupuesto")
$wsPriorizacion = $wb.Worksheets.Item("Priorización")
$wsHoja1        = $wb.Worksheets.Item("Hoja1")

# --- Shift the helper list down one row and insert the new 0 .. 1 range ---
for ($r = 1; $r -le 11; $r++) {
    $wsHoja1.Cells.Item($r, 1).Value = ($r - 1) / 10.0
}
# Make sure the newly-added row keeps the same "0.0" number format as the
# rest of the column (it has no predecessor to inherit the style from).
$wsHoja1.Cells.Item(11, 1).NumberFormat = "0.0"

# --- Make "Priorización" the active sheet/tab (was "Presupuesto") ---
$wsPriorizacion.Activate()
